## Emailing bacs creating month report and sending month end report code implementation done
# Populate two new SDLT entries (rows 3 and 4) and refresh the first data
# row (row 2) with the current month's case - "New Build Purchase" records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : NNT1028 - Beenish Ali -----------------------------------
$ws.Range("A2").Value = "NNT1028"
$ws.Range("B2").Value = "New Build Purchase"
$ws.Range("D2").Value = "'False"
$ws.Range("E2").Value = "Beenish Ali"
$ws.Range("I2").Value = "Plot 7.08 Belgrave Village Birmingham             "
$ws.Range("K2").Value = "Plot 7.08 Belgrave Village Birmingham               "
$ws.Range("L2").Value = "£187,000.00"
$ws.Range("M2").Value = "16 September 2022"
$ws.Range("N2").Value = "19 July 2023"
$ws.Range("R2").Value = "Howard Kennedy"
$ws.Range("S2").Value = "1 London Bridge, SE1 9BG"
$ws.Range("T2").Value = "Plot 7.08 Belgrave Village Birmingham   "
$ws.Range("V2").Value = "Individuals"
$ws.Range("AC2").Value = "Beenish Ali"
$ws.Range("AD2").Value = "SC809688A"

# --- Row 3 : NBT1872 - Marius Catalin Voica and Maria Magdalena Voica -
$ws.Range("A3").Value = "NBT1872"
$ws.Range("B3").Value = "New Build Purchase"
$ws.Range("D3").Value = "'False"
$ws.Range("E3").Value = "Marius Catalin Voica and Maria Magdalena Voica"
$ws.Range("I3").Value = " Plot 274  Wintringham Park, St. Neots          "
$ws.Range("K3").Value = "  Plot 274  Wintringham Park, St. Neots           "
$ws.Range("L3").Value = "£259,950.00"
$ws.Range("M3").Value = "20 July 2023"
$ws.Range("N3").Value = "20 July 2023"
$ws.Range("R3").Value = "Birketts"
$ws.Range("S3").Value = "Kingfisher House, 1 Gilders Way, Norwich, Norfolk, NR3 1UB"
$ws.Range("T3").Value = " Plot 274 Wintringham Park, St. Neots "
$ws.Range("V3").Value = "Individuals"
$ws.Range("AC3").Value = "Marius Catalin Voica"
$ws.Range("AD3").Value = "SS438175D"

# --- Row 4 : NNT1227 - Silje Merete Sathren Gronning -------------------
$ws.Range("A4").Value = "NNT1227"
$ws.Range("B4").Value = "New Build Purchase"
$ws.Range("D4").Value = "'False"
$ws.Range("E4").Value = "Silje Merete Sathren Gronning"
$ws.Range("I4").Value = "Flat 50, 8 Newton Street, London             "
$ws.Range("K4").Value = "  Plot 200  Acton Gardens            "
$ws.Range("L4").Value = "£460,000.00"
$ws.Range("M4").Value = "12 July 2023"
$ws.Range("N4").Value = "20 July 2023"
$ws.Range("R4").Value = "Countryside Properties (UK) Limited"
$ws.Range("S4").Value = "DX 124280 , BRENTWOOD 4"
$ws.Range("T4").Value = " Plot 200 Acton Gardens  "
$ws.Range("V4").Value = "Individuals"
$ws.Range("AC4").Value = "Silje Merete Sathren Gronning"
$ws.Range("AD4").Value = "SY217118C"
